$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 80, shifting existing rows 80-154 down to 82-156
$ws.Rows("80:81").Insert()

# Row 80 - new weekly data entry (Lane Late / Primera)
$ws.Cells.Item(80,1).Value2  = 11
$ws.Cells.Item(80,2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(80,3).Value2  = "Bíobío"
$ws.Cells.Item(80,4).Value2  = 44484
$ws.Cells.Item(80,5).Value2  = 8
$ws.Cells.Item(80,6).Value2  = "Fruta"
$ws.Cells.Item(80,7).Value2  = 100102
$ws.Cells.Item(80,8).Value2  = "Cítricos"
$ws.Cells.Item(80,9).Value2  = 100102005
$ws.Cells.Item(80,10).Value2 = "Naranja"
$ws.Cells.Item(80,11).Value2 = "Lane Late"
$ws.Cells.Item(80,12).Value2 = "Primera"
$ws.Cells.Item(80,13).Value2 = 100
$ws.Cells.Item(80,14).Value2 = 8000
$ws.Cells.Item(80,15).Value2 = 8000
$ws.Cells.Item(80,16).Value2 = 8000
$ws.Cells.Item(80,17).Value2 = "$/caja 15 kilos empedrada"
$ws.Cells.Item(80,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(80,19).Value2 = 533
$ws.Cells.Item(80,20).Value2 = 15

# Row 81 - new weekly data entry (Lane Late / Segunda)
$ws.Cells.Item(81,1).Value2  = 11
$ws.Cells.Item(81,2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(81,3).Value2  = "Bíobío"
$ws.Cells.Item(81,4).Value2  = 44484
$ws.Cells.Item(81,5).Value2  = 8
$ws.Cells.Item(81,6).Value2  = "Fruta"
$ws.Cells.Item(81,7).Value2  = 100102
$ws.Cells.Item(81,8).Value2  = "Cítricos"
$ws.Cells.Item(81,9).Value2  = 100102005
$ws.Cells.Item(81,10).Value2 = "Naranja"
$ws.Cells.Item(81,11).Value2 = "Lane Late"
$ws.Cells.Item(81,12).Value2 = "Segunda"
$ws.Cells.Item(81,13).Value2 = 100
$ws.Cells.Item(81,14).Value2 = 7000
$ws.Cells.Item(81,15).Value2 = 7000
$ws.Cells.Item(81,16).Value2 = 7000
$ws.Cells.Item(81,17).Value2 = "$/caja 15 kilos empedrada"
$ws.Cells.Item(81,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(81,19).Value2 = 467
$ws.Cells.Item(81,20).Value2 = 15
